$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Save" header in H1, reusing the exact formatting of the other
# headers (e.g. G1) by copying its format over.
$ws.Range("G1").Copy() | Out-Null
$ws.Range("H1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("H1").Value = "Save"

# Add the value for the new Save column in H2
$ws.Range("H2").Value = 1
